$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "score on submission" for the unet_augm_diceloss row
$ws.Range("D4").Value = 0.725

# Append a new results row for the new submission (unet pc augm diceloss, score 0.897)
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats, reuse the date style
$ws.Range("A6").Value = 44449
$ws.Range("B6").Value = "model_floodwater_unet_pc_augm_diceloss"
$ws.Range("C6").Value = 0.715
$ws.Range("D6").Value = 0.897
$ws.Range("E6").Value = "['kuo', 'wvy', 'awc']"

# Move active selection like Excel would after entering the new row
$ws.Range("D7").Select()
